$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-18 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-19 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("70×90=", $true, $false, $false, $false, $false, $true, 1, $false, "14×61=", 2) | Out-Null
$d.Content.Find.Execute("69×82=", $true, $false, $false, $false, $false, $true, 1, $false, "38×79=", 2) | Out-Null
$d.Content.Find.Execute("98×92=", $true, $false, $false, $false, $false, $true, 1, $false, "82×53=", 2) | Out-Null
$d.Content.Find.Execute("97×91=", $true, $false, $false, $false, $false, $true, 1, $false, "36×59=", 2) | Out-Null
$d.Content.Find.Execute("55×20=", $true, $false, $false, $false, $false, $true, 1, $false, "55×36=", 2) | Out-Null
$d.Content.Find.Execute("37×90=", $true, $false, $false, $false, $false, $true, 1, $false, "52×44=", 2) | Out-Null
$d.Content.Find.Execute("27×86=", $true, $false, $false, $false, $false, $true, 1, $false, "66×85=", 2) | Out-Null
$d.Content.Find.Execute("40×95=", $true, $false, $false, $false, $false, $true, 1, $false, "98×41=", 2) | Out-Null
$d.Content.Find.Execute("59×38=", $true, $false, $false, $false, $false, $true, 1, $false, "43×50=", 2) | Out-Null
$d.Content.Find.Execute("45×98=", $true, $false, $false, $false, $false, $true, 1, $false, "96×80=", 2) | Out-Null
$d.Content.Find.Execute("34×25=", $true, $false, $false, $false, $false, $true, 1, $false, "71×45=", 2) | Out-Null
$d.Content.Find.Execute("20×73=", $true, $false, $false, $false, $false, $true, 1, $false, "94×21=", 2) | Out-Null
$d.Content.Find.Execute("12×98=", $true, $false, $false, $false, $false, $true, 1, $false, "37×37=", 2) | Out-Null
$d.Content.Find.Execute("81×27=", $true, $false, $false, $false, $false, $true, 1, $false, "55×28=", 2) | Out-Null
$d.Content.Find.Execute("99×16=", $true, $false, $false, $false, $false, $true, 1, $false, "29×67=", 2) | Out-Null
$d.Content.Find.Execute("81×35=", $true, $false, $false, $false, $false, $true, 1, $false, "98×53=", 2) | Out-Null
$d.Content.Find.Execute("66×54=", $true, $false, $false, $false, $false, $true, 1, $false, "44×46=", 2) | Out-Null
$d.Content.Find.Execute("25×50=", $true, $false, $false, $false, $false, $true, 1, $false, "47×38=", 2) | Out-Null
$d.Content.Find.Execute("82×84=", $true, $false, $false, $false, $false, $true, 1, $false, "17×86=", 2) | Out-Null
$d.Content.Find.Execute("93×28=", $true, $false, $false, $false, $false, $true, 1, $false, "80×85=", 2) | Out-Null
$d.Content.Find.Execute("70×37=", $true, $false, $false, $false, $false, $true, 1, $false, "78×94=", 2) | Out-Null
$d.Content.Find.Execute("13×70=", $true, $false, $false, $false, $false, $true, 1, $false, "44×14=", 2) | Out-Null
$d.Content.Find.Execute("75×54=", $true, $false, $false, $false, $false, $true, 1, $false, "66×63=", 2) | Out-Null
$d.Content.Find.Execute("29×82=", $true, $false, $false, $false, $false, $true, 1, $false, "98×18=", 2) | Out-Null
$d.Content.Find.Execute("25×26=", $true, $false, $false, $false, $false, $true, 1, $false, "29×27=", 2) | Out-Null

Write-Output "Done applying replacements."
